# Auto-generated edit script applying diff changes to 上海-漫展信息.xlsx
$wb = $excel.ActiveWorkbook

# --- sheet1 (展览) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value2 = 355
$ws.Cells.Item(3, 6).Value2 = 1271
$ws.Cells.Item(4, 6).Value2 = 84
$ws.Cells.Item(5, 6).Value2 = 2055
$ws.Cells.Item(6, 6).Value2 = 106
$ws.Cells.Item(7, 6).Value2 = 798
$ws.Cells.Item(8, 6).Value2 = 15
$ws.Cells.Item(12, 6).Value2 = 1061
$ws.Cells.Item(13, 6).Value2 = 793
$ws.Cells.Item(15, 6).Value2 = 649
$ws.Cells.Item(16, 6).Value2 = 1262
$ws.Cells.Item(20, 6).Value2 = 713
$ws.Cells.Item(21, 6).Value2 = 77
$ws.Cells.Item(23, 6).Value2 = 91
$ws.Cells.Item(25, 6).Value2 = 1205
$ws.Cells.Item(26, 6).Value2 = 122
$ws.Cells.Item(27, 6).Value2 = 417
$ws.Cells.Item(29, 6).Value2 = 4983
$ws.Cells.Item(30, 6).Value2 = 228
$ws.Cells.Item(33, 6).Value2 = 5766
$ws.Cells.Item(33, 9).Value2 = "//i1.hdslb.com/bfs/openplatform/202404/UH8lEMFK1714277926888.jpeg"
$ws.Cells.Item(34, 6).Value2 = 959
$ws.Cells.Item(35, 6).Value2 = 576
$ws.Cells.Item(38, 6).Value2 = 1035
$ws.Cells.Item(41, 6).Value2 = 650

# --- sheet2 (演出) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(5, 6).Value2 = 2078
$ws.Cells.Item(8, 6).Value2 = 110
$ws.Cells.Item(15, 6).Value2 = 646
$ws.Cells.Item(16, 6).Value2 = 646
$ws.Cells.Item(23, 6).Value2 = 16
$ws.Cells.Item(24, 6).Value2 = 97
$ws.Cells.Item(28, 6).Value2 = 1708
$ws.Cells.Item(34, 6).Value2 = 95
$ws.Cells.Item(37, 6).Value2 = 60
$ws.Cells.Item(41, 6).Value2 = 888
$ws.Cells.Item(42, 6).Value2 = 476
$ws.Cells.Item(47, 6).Value2 = 7

# --- sheet3 (本地生活) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(7, 6).Value2 = 349
$ws.Cells.Item(8, 6).Value2 = 201

# --- sheet4 (全部类型) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value2 = 355
$ws.Cells.Item(6, 6).Value2 = 1271
$ws.Cells.Item(7, 6).Value2 = 84
$ws.Cells.Item(8, 6).Value2 = 349
$ws.Cells.Item(9, 6).Value2 = 201
$ws.Cells.Item(10, 6).Value2 = 201
$ws.Cells.Item(11, 6).Value2 = 2078
$ws.Cells.Item(12, 6).Value2 = 2055
$ws.Cells.Item(14, 6).Value2 = 798
$ws.Cells.Item(15, 6).Value2 = 110
$ws.Cells.Item(16, 6).Value2 = 15
$ws.Cells.Item(19, 6).Value2 = 1061
$ws.Cells.Item(20, 6).Value2 = 793
$ws.Cells.Item(23, 6).Value2 = 649
$ws.Cells.Item(24, 6).Value2 = 1262
$ws.Cells.Item(27, 6).Value2 = 713
$ws.Cells.Item(28, 6).Value2 = 77
$ws.Cells.Item(30, 6).Value2 = 646
$ws.Cells.Item(32, 6).Value2 = 1205
$ws.Cells.Item(33, 6).Value2 = 122
$ws.Cells.Item(35, 6).Value2 = 417
$ws.Cells.Item(36, 6).Value2 = 4983
$ws.Cells.Item(37, 6).Value2 = 228
$ws.Cells.Item(40, 6).Value2 = 5766
$ws.Cells.Item(40, 9).Value2 = "//i1.hdslb.com/bfs/openplatform/202404/UH8lEMFK1714277926888.jpeg"
$ws.Cells.Item(41, 6).Value2 = 959
$ws.Cells.Item(42, 6).Value2 = 1708
$ws.Cells.Item(43, 6).Value2 = 576
$ws.Cells.Item(45, 6).Value2 = 1035
$ws.Cells.Item(46, 6).Value2 = 650
$ws.Cells.Item(47, 6).Value2 = 60
$ws.Cells.Item(50, 6).Value2 = 888
$ws.Cells.Item(51, 6).Value2 = 476
